# Horarios actualizados Linea 141 - 1226
#
# Periodic refresh of the scraped live bus-arrival data (source: data/
# horarios-141-2026-01-31.xlsx). The scraper re-ran at 07:21:42, so on every
# sheet the "Ultima actualizacion" / "Total filas" header is bumped and the
# tail of the arrivals table (rows whose Hora_Scrap was still the previous
# 06:5x:xx/06:xx:xx scrape, i.e. not-yet-arrived buses) is rewritten with the
# newer figures; several brand-new arrivals were appended at the bottom.
# Earlier rows (already-arrived buses, stable since an older scrape) are
# untouched.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "LP1912"  ->  used range grows to A1:E73
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Header metadata: last-scrape timestamp + row count
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:21:42'
$ws.Cells.Item(3,1).Value = 'Total filas: 68'

# Rows 43-65: updated in place (later scrape picked up new Hora_Scrap/Hora_Llegada/Minutos, occasionally a
# different stop, for arrivals still in the future)
$ws.Cells.Item(43,1).Value = '07:21:42'
$ws.Cells.Item(43,2).Value = '07:21'
$ws.Cells.Item(43,3).Value = '215A_EL PATO'
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,5).Value = 'LP1912'
$ws.Cells.Item(44,1).Value = '07:21:42'
$ws.Cells.Item(44,2).Value = '07:23'
$ws.Cells.Item(44,3).Value = '16_SANTA ANA'
$ws.Cells.Item(44,4).Value = 2
$ws.Cells.Item(44,5).Value = 'LP1912'
$ws.Cells.Item(45,1).Value = '06:58:01'
$ws.Cells.Item(45,2).Value = '07:24'
$ws.Cells.Item(45,3).Value = '16_SANTA ANA'
$ws.Cells.Item(45,4).Value = 26
$ws.Cells.Item(45,5).Value = 'LP1912'
$ws.Cells.Item(46,1).Value = '07:21:42'
$ws.Cells.Item(46,2).Value = '07:29'
$ws.Cells.Item(46,3).Value = '14_ABASTO'
$ws.Cells.Item(46,4).Value = 8
$ws.Cells.Item(46,5).Value = 'LP1912'
$ws.Cells.Item(47,1).Value = '07:21:42'
$ws.Cells.Item(47,2).Value = '07:33'
$ws.Cells.Item(47,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(47,4).Value = 12
$ws.Cells.Item(47,5).Value = 'LP1912'
$ws.Cells.Item(48,1).Value = '06:58:01'
$ws.Cells.Item(48,2).Value = '07:34'
$ws.Cells.Item(48,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(48,4).Value = 36
$ws.Cells.Item(48,5).Value = 'LP1912'
$ws.Cells.Item(49,1).Value = '07:21:42'
$ws.Cells.Item(49,2).Value = '07:36'
$ws.Cells.Item(49,3).Value = '17X38_ROMERO'
$ws.Cells.Item(49,4).Value = 15
$ws.Cells.Item(49,5).Value = 'LP1912'
$ws.Cells.Item(50,1).Value = '07:21:42'
$ws.Cells.Item(50,2).Value = '07:36'
$ws.Cells.Item(50,3).Value = '27_EL RETIRO'
$ws.Cells.Item(50,4).Value = 15
$ws.Cells.Item(50,5).Value = 'LP1912'
$ws.Cells.Item(51,1).Value = '06:58:01'
$ws.Cells.Item(51,2).Value = '07:37'
$ws.Cells.Item(51,3).Value = '27_EL RETIRO'
$ws.Cells.Item(51,4).Value = 39
$ws.Cells.Item(51,5).Value = 'LP1912'
$ws.Cells.Item(52,1).Value = '07:21:42'
$ws.Cells.Item(52,2).Value = '07:41'
$ws.Cells.Item(52,3).Value = '16_SANTA ANA'
$ws.Cells.Item(52,4).Value = 20
$ws.Cells.Item(52,5).Value = 'LP1912'
$ws.Cells.Item(53,1).Value = '07:21:42'
$ws.Cells.Item(53,2).Value = '07:43'
$ws.Cells.Item(53,3).Value = '10_OLMOS'
$ws.Cells.Item(53,4).Value = 22
$ws.Cells.Item(53,5).Value = 'LP1912'
$ws.Cells.Item(54,1).Value = '06:58:01'
$ws.Cells.Item(54,2).Value = '07:44'
$ws.Cells.Item(54,3).Value = '10_OLMOS'
$ws.Cells.Item(54,4).Value = 46
$ws.Cells.Item(54,5).Value = 'LP1912'
$ws.Cells.Item(55,1).Value = '07:21:42'
$ws.Cells.Item(55,2).Value = '07:49'
$ws.Cells.Item(55,3).Value = '15_ABASTO'
$ws.Cells.Item(55,4).Value = 28
$ws.Cells.Item(55,5).Value = 'LP1912'
$ws.Cells.Item(56,1).Value = '07:21:42'
$ws.Cells.Item(56,2).Value = '07:58'
$ws.Cells.Item(56,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(56,4).Value = 37
$ws.Cells.Item(56,5).Value = 'LP1912'
$ws.Cells.Item(57,1).Value = '07:21:42'
$ws.Cells.Item(57,2).Value = '07:59'
$ws.Cells.Item(57,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(57,4).Value = 38
$ws.Cells.Item(57,5).Value = 'LP1912'
$ws.Cells.Item(58,1).Value = '06:58:01'
$ws.Cells.Item(58,2).Value = '08:00'
$ws.Cells.Item(58,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(58,4).Value = 62
$ws.Cells.Item(58,5).Value = 'LP1912'
$ws.Cells.Item(59,1).Value = '06:46:06'
$ws.Cells.Item(59,2).Value = '08:03'
$ws.Cells.Item(59,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(59,4).Value = 77
$ws.Cells.Item(59,5).Value = 'LP1912'
$ws.Cells.Item(60,1).Value = '07:21:42'
$ws.Cells.Item(60,2).Value = '08:03'
$ws.Cells.Item(60,3).Value = '17X38_ROMERO'
$ws.Cells.Item(60,4).Value = 42
$ws.Cells.Item(60,5).Value = 'LP1912'
$ws.Cells.Item(61,1).Value = '06:58:01'
$ws.Cells.Item(61,2).Value = '08:04'
$ws.Cells.Item(61,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(61,4).Value = 66
$ws.Cells.Item(61,5).Value = 'LP1912'
$ws.Cells.Item(62,1).Value = '07:21:42'
$ws.Cells.Item(62,2).Value = '08:14'
$ws.Cells.Item(62,3).Value = '10_OLMOS'
$ws.Cells.Item(62,4).Value = 53
$ws.Cells.Item(62,5).Value = 'LP1912'
$ws.Cells.Item(63,1).Value = '07:21:42'
$ws.Cells.Item(63,2).Value = '08:19'
$ws.Cells.Item(63,3).Value = '15_ABASTO'
$ws.Cells.Item(63,4).Value = 58
$ws.Cells.Item(63,5).Value = 'LP1912'
$ws.Cells.Item(64,1).Value = '07:21:42'
$ws.Cells.Item(64,2).Value = '08:29'
$ws.Cells.Item(64,3).Value = '14_ABASTO'
$ws.Cells.Item(64,4).Value = 68
$ws.Cells.Item(64,5).Value = 'LP1912'
$ws.Cells.Item(65,1).Value = '06:58:01'
$ws.Cells.Item(65,2).Value = '08:30'
$ws.Cells.Item(65,3).Value = '14_ABASTO'
$ws.Cells.Item(65,4).Value = 92
$ws.Cells.Item(65,5).Value = 'LP1912'

# Rows 66-73: brand-new arrivals appended by the scraper
$ws.Cells.Item(66,1).Value = '07:21:42'
$ws.Cells.Item(66,2).Value = '08:33'
$ws.Cells.Item(66,3).Value = '215C_EL PATO'
$ws.Cells.Item(66,4).Value = 72
$ws.Cells.Item(66,5).Value = 'LP1912'
$ws.Cells.Item(67,1).Value = '06:58:01'
$ws.Cells.Item(67,2).Value = '08:34'
$ws.Cells.Item(67,3).Value = '215C_EL PATO'
$ws.Cells.Item(67,4).Value = 96
$ws.Cells.Item(67,5).Value = 'LP1912'
$ws.Cells.Item(68,1).Value = '07:21:42'
$ws.Cells.Item(68,2).Value = '08:48'
$ws.Cells.Item(68,3).Value = '215A_EL PATO'
$ws.Cells.Item(68,4).Value = 87
$ws.Cells.Item(68,5).Value = 'LP1912'
$ws.Cells.Item(69,1).Value = '07:21:42'
$ws.Cells.Item(69,2).Value = '08:51'
$ws.Cells.Item(69,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(69,4).Value = 90
$ws.Cells.Item(69,5).Value = 'LP1912'
$ws.Cells.Item(70,1).Value = '07:21:42'
$ws.Cells.Item(70,2).Value = '08:59'
$ws.Cells.Item(70,3).Value = '215B_EL PATO'
$ws.Cells.Item(70,4).Value = 98
$ws.Cells.Item(70,5).Value = 'LP1912'
$ws.Cells.Item(71,1).Value = '07:21:42'
$ws.Cells.Item(71,2).Value = '09:14'
$ws.Cells.Item(71,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(71,4).Value = 113
$ws.Cells.Item(71,5).Value = 'LP1912'
$ws.Cells.Item(72,1).Value = '07:21:42'
$ws.Cells.Item(72,2).Value = '09:16'
$ws.Cells.Item(72,3).Value = '27_EL RETIRO'
$ws.Cells.Item(72,4).Value = 115
$ws.Cells.Item(72,5).Value = 'LP1912'
$ws.Cells.Item(73,1).Value = '07:21:42'
$ws.Cells.Item(73,2).Value = '09:18'
$ws.Cells.Item(73,3).Value = '215_EL PELIGRO'
$ws.Cells.Item(73,4).Value = 117
$ws.Cells.Item(73,5).Value = 'LP1912'

# ---------------------------------------------------------------
# Sheet 2: "LP1912-215"  ->  used range grows to A1:E20
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Header metadata: last-scrape timestamp + row count
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:21:42'
$ws.Cells.Item(3,1).Value = 'Total filas: 15'

# Rows 15-17: updated in place (later scrape picked up new Hora_Scrap/Hora_Llegada/Minutos, occasionally a
# different stop, for arrivals still in the future)
$ws.Cells.Item(15,1).Value = '07:21:42'
$ws.Cells.Item(15,2).Value = '07:21'
$ws.Cells.Item(15,3).Value = '215A_EL PATO'
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 'LP1912'
$ws.Cells.Item(16,1).Value = '07:21:42'
$ws.Cells.Item(16,2).Value = '08:33'
$ws.Cells.Item(16,3).Value = '215C_EL PATO'
$ws.Cells.Item(16,4).Value = 72
$ws.Cells.Item(16,5).Value = 'LP1912'
$ws.Cells.Item(17,1).Value = '06:58:01'
$ws.Cells.Item(17,2).Value = '08:34'
$ws.Cells.Item(17,3).Value = '215C_EL PATO'
$ws.Cells.Item(17,4).Value = 96
$ws.Cells.Item(17,5).Value = 'LP1912'

# Rows 18-20: brand-new arrivals appended by the scraper
$ws.Cells.Item(18,1).Value = '07:21:42'
$ws.Cells.Item(18,2).Value = '08:48'
$ws.Cells.Item(18,3).Value = '215A_EL PATO'
$ws.Cells.Item(18,4).Value = 87
$ws.Cells.Item(18,5).Value = 'LP1912'
$ws.Cells.Item(19,1).Value = '07:21:42'
$ws.Cells.Item(19,2).Value = '08:59'
$ws.Cells.Item(19,3).Value = '215B_EL PATO'
$ws.Cells.Item(19,4).Value = 98
$ws.Cells.Item(19,5).Value = 'LP1912'
$ws.Cells.Item(20,1).Value = '07:21:42'
$ws.Cells.Item(20,2).Value = '09:18'
$ws.Cells.Item(20,3).Value = '215_EL PELIGRO'
$ws.Cells.Item(20,4).Value = 117
$ws.Cells.Item(20,5).Value = 'LP1912'

# ---------------------------------------------------------------
# Sheet 3: "6203-6173"  ->  used range grows to A1:E16
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

# Header metadata: last-scrape timestamp + row count
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:21:42'
$ws.Cells.Item(3,1).Value = 'Total filas: 11'

# Rows 11-14: updated in place (later scrape picked up new Hora_Scrap/Hora_Llegada/Minutos, occasionally a
# different stop, for arrivals still in the future)
$ws.Cells.Item(11,1).Value = '07:21:42'
$ws.Cells.Item(11,2).Value = '07:35'
$ws.Cells.Item(11,3).Value = '215A_LA PLATA'
$ws.Cells.Item(11,4).Value = 14
$ws.Cells.Item(11,5).Value = 'L6173'
$ws.Cells.Item(12,1).Value = '06:58:01'
$ws.Cells.Item(12,2).Value = '07:36'
$ws.Cells.Item(12,3).Value = '215A_LA PLATA'
$ws.Cells.Item(12,4).Value = 38
$ws.Cells.Item(12,5).Value = 'L6173'
$ws.Cells.Item(13,1).Value = '07:21:42'
$ws.Cells.Item(13,2).Value = '08:09'
$ws.Cells.Item(13,3).Value = '215A_LA PLATA'
$ws.Cells.Item(13,4).Value = 48
$ws.Cells.Item(13,5).Value = 'L6173'
$ws.Cells.Item(14,1).Value = '06:46:06'
$ws.Cells.Item(14,2).Value = '08:10'
$ws.Cells.Item(14,3).Value = '215A_LA PLATA'
$ws.Cells.Item(14,4).Value = 84
$ws.Cells.Item(14,5).Value = 'L6173'

# Rows 15-16: brand-new arrivals appended by the scraper
$ws.Cells.Item(15,1).Value = '07:21:42'
$ws.Cells.Item(15,2).Value = '08:23'
$ws.Cells.Item(15,3).Value = '215C_LA PLATA'
$ws.Cells.Item(15,4).Value = 62
$ws.Cells.Item(15,5).Value = 'L6203'
$ws.Cells.Item(16,1).Value = '06:58:01'
$ws.Cells.Item(16,2).Value = '08:52'
$ws.Cells.Item(16,3).Value = '215A_LA PLATA'
$ws.Cells.Item(16,4).Value = 114
$ws.Cells.Item(16,5).Value = 'L6173'

